# AMMR 3.0 Thoracic model development - updated Wilke spinal disc pressure
# validation data (new AMS force measurements in column C). The dependent
# ratios in column D (=C/C5) and the linked chart series recompute
# automatically from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Updated AMS [N] measurements (column C, rows 2-9)
$ws.Range("C2").Value = 151.9
$ws.Range("C3").Value = 287.9
$ws.Range("C4").Value = 457.6
$ws.Range("C5").Value = 595.70000000000005
$ws.Range("C6").Value = 1282.8
$ws.Range("C7").Value = 1181.5
$ws.Range("C8").Value = 2813.4
$ws.Range("C9").Value = 2117.8000000000002

# Reflect the author's final cursor position on the sheet
[void]$ws.Range("G14").Select()
